$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.052.43"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.911.67"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8279"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3230"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07030"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08036"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7513"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "1.909.82"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.228"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "30.049.41"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.937"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007776"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "2.156.30"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.992"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1621"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +24.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.266"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.086"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.369"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.519"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.309"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05601"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.099"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7366"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01920"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.796"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4442"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8424"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.620"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.772"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "983.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.06%  "
$ws.Range("D50").Value = "2.062.77"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
